# TC10/TC11 Kaman testcase workbook update
# - Rename the first sheet (TC11 -> TC10)
# - Update the testcase-name cell to match
# - Insert a WAIT step after the "press enter" step
# - Insert a TINY_SCROLL_DOWN step after the "verify search header" step
# - Update Testdata sheet: SearchBoxHomePage becomes numeric 200 (was text "200CL"),
#   SearchHeader text drops the "CL" suffix, and two new EleType rows are appended

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: rename ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "TC10_Verify_PDP_From_Search"

# Testcase name cell
$ws1.Range("A2").Value = "TC10_Verify_PDP_From_Search_results"

# Insert the new "WAIT" row right after the PRESS_ENTER step (old row 6, which
# held VERIFY_WEBELEMENT_PRESENT/SearchResults, shifts down to make room)
$ws1.Rows.Item(6).Insert()
$ws1.Range("B6").Value = "WAIT"

# Insert the new "TINY_SCROLL_DOWN" row right after the VERIFY_TEXT_PRESENT /
# SearchHeader step (now at row 8 after the previous insert)
$ws1.Rows.Item(9).Insert()
$ws1.Range("B9").Value = "TINY_SCROLL_DOWN"

# ---- Sheet 2: Testdata ----
$ws2 = $wb.Worksheets.Item(2)

# SearchBoxHomePage value: "200CL" (text) -> 200 (number)
$ws2.Range("B3").Value = 200

# SearchHeader expected text drops the CL suffix
$ws2.Range("B5").Value = "Showing Results for ""200"""

# New rows: EleType1 / EleType2, both JSElement
$ws2.Range("A7").Value = "EleType1"
$ws2.Range("B7").Value = "JSElement"
$ws2.Range("A8").Value = "EleType2"
$ws2.Range("B8").Value = "JSElement"

# Sheet2 selection as left by the author
$ws2.Range("B23").Select()

# Re-activate sheet1 last so it stays the visible/tabSelected sheet, and set
# its own lingering selection
$ws1.Activate()
$ws1.Range("A3:XFD7").Select()
